$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A10").Value = "'009"
$ws.Range("B10").Value = "el 006 con tendencia"
$ws.Range("B10").Select() | Out-Null
